$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 33
$ws.Range("H33").Value = 834.1111
$ws.Range("I33").Value = 539.8333
$ws.Range("J33").Value = 1422.6666
$ws.Range("K33").Value = 539.8333
$ws.Range("L33").Value = 1422.6666
$ws.Range("M33").Value = -310.8333
$ws.Range("N33").Value = -1880.6666

# ALC row 138
$ws.Range("H138").Value = 4399.3047
$ws.Range("I138").Value = 2736.8708
$ws.Range("J138").Value = 5409.8037
$ws.Range("K138").Value = 8210.6124
$ws.Range("L138").Value = 16229.4111
$ws.Range("M138").Value = -3070.6124
$ws.Range("N138").Value = -26509.4111

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 17135.307
$ws.Range("I32").Value = 14209.213
$ws.Range("K32").Value = 14209.213
$ws.Range("M32").Value = -13922.213

# ARM row 45
$ws.Range("H45").Value = 1391.6111
$ws.Range("I45").Value = 1090.7059
$ws.Range("K45").Value = 1090.7059
$ws.Range("M45").Value = -713.7058999999999

# ARM row 74
$ws.Range("H74").Value = 1822.3715
$ws.Range("I74").Value = 1392.069
$ws.Range("K74").Value = 1392.069
$ws.Range("M74").Value = -518.069

# ARM row 77
$ws.Range("H77").Value = 1822.3715
$ws.Range("I77").Value = 1392.069
$ws.Range("K77").Value = 6960.344999999999
$ws.Range("M77").Value = -2592.344999999999

# ARM row 122
$ws.Range("H122").Value = 2963.5715
$ws.Range("I122").Value = 2021.4
$ws.Range("J122").Value = 5319
$ws.Range("K122").Value = 6064.200000000001
$ws.Range("L122").Value = 15957
$ws.Range("M122").Value = -3614.200000000001
$ws.Range("N122").Value = -20857

# ARM row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ARM row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# BSM row 51
$ws.Range("H51").Value = 29780
$ws.Range("J51").Value = 29780
$ws.Range("L51").Value = 29780
$ws.Range("N51").Value = -30762

# BSM row 105
$ws.Range("H105").Value = 2352.1
$ws.Range("I105").Value = 2102.0715
$ws.Range("K105").Value = 2102.0715
$ws.Range("M105").Value = -355.0715

# BSM row 107
$ws.Range("H107").Value = 3292.1538
$ws.Range("I107").Value = 2644.2222
$ws.Range("K107").Value = 2644.2222
$ws.Range("M107").Value = -724.2222000000002

$ws = $wb.Worksheets.Item("CRP")
# CRP row 22
$ws.Range("H22").Value = 949.5
$ws.Range("I22").Value = 349.35715
$ws.Range("J22").Value = 3050
$ws.Range("K22").Value = 349.35715
$ws.Range("L22").Value = 3050
$ws.Range("M22").Value = 0.6428500000000099
$ws.Range("N22").Value = -3750

# CRP row 31
$ws.Range("H31").Value = 1495447.4
$ws.Range("I31").Value = 2042528.9
$ws.Range("J31").Value = 6170.1113
$ws.Range("K31").Value = 2042528.9
$ws.Range("L31").Value = 6170.1113
$ws.Range("M31").Value = -2042233.9
$ws.Range("N31").Value = -6760.1113

# CRP row 34
$ws.Range("H34").Value = 1495447.4
$ws.Range("I34").Value = 2042528.9
$ws.Range("J34").Value = 6170.1113
$ws.Range("K34").Value = 2042528.9
$ws.Range("L34").Value = 6170.1113
$ws.Range("M34").Value = -2042326.9
$ws.Range("N34").Value = -6574.1113

# CRP row 123
$ws.Range("H123").Value = 30956
$ws.Range("J123").Value = 30956
$ws.Range("L123").Value = 30956
$ws.Range("N123").Value = -40756

# CRP row 141
$ws.Range("H141").Value = 30500
$ws.Range("J141").Value = 32109.375
$ws.Range("L141").Value = 32109.375
$ws.Range("N141").Value = -42469.375

$ws = $wb.Worksheets.Item("CUL")
# CUL row 87
$ws.Range("H87").Value = 7467
$ws.Range("I87").Value = 3213
$ws.Range("K87").Value = 9639
$ws.Range("M87").Value = -8391

# CUL row 90
$ws.Range("H90").Value = 7467
$ws.Range("I90").Value = 3213
$ws.Range("K90").Value = 28917
$ws.Range("M90").Value = -22677

# CUL row 113
$ws.Range("H113").Value = 2778629.8
$ws.Range("J113").Value = 930.2759
$ws.Range("L113").Value = 2790.8277
$ws.Range("N113").Value = -7130.8277

# CUL row 140
$ws.Range("H140").Value = 5610.222
$ws.Range("I140").Value = 2100
$ws.Range("J140").Value = 9998
$ws.Range("K140").Value = 6300
$ws.Range("L140").Value = 29994
$ws.Range("M140").Value = -1120
$ws.Range("N140").Value = -40354

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Range("H2").Value = 88
$ws.Range("I2").Value = 100.90909
$ws.Range("K2").Value = 100.90909
$ws.Range("M2").Value = 12.09090999999999

# GSM row 122
$ws.Range("H122").Value = 5113
$ws.Range("I122").Value = 3416.6667
$ws.Range("J122").Value = 10202
$ws.Range("K122").Value = 10250.0001
$ws.Range("L122").Value = 30606
$ws.Range("M122").Value = -7800.000100000001
$ws.Range("N122").Value = -35506

# GSM row 126
$ws.Range("H126").Value = 671078.7
$ws.Range("I126").Value = 3374.6667
$ws.Range("J126").Value = 838004.7
$ws.Range("K126").Value = 10124.0001
$ws.Range("L126").Value = 2514014.1
$ws.Range("M126").Value = -7654.000100000001
$ws.Range("N126").Value = -2518954.1

# GSM row 132
$ws.Range("H132").Value = 3720.1702
$ws.Range("I132").Value = 3848.276
$ws.Range("J132").Value = 3513.7778
$ws.Range("K132").Value = 11544.828
$ws.Range("L132").Value = 10541.3334
$ws.Range("M132").Value = -9014.828
$ws.Range("N132").Value = -15601.3334

$ws = $wb.Worksheets.Item("LTW")
# LTW row 24
$ws.Range("H24").Value = 19998.666
$ws.Range("J24").Value = 19998.666
$ws.Range("L24").Value = 19998.666
$ws.Range("N24").Value = -20684.666

# LTW row 82
$ws.Range("H82").Value = 2424.8235
$ws.Range("I82").Value = 1863.1428
$ws.Range("J82").Value = 2818
$ws.Range("K82").Value = 1863.1428
$ws.Range("L82").Value = 2818
$ws.Range("M82").Value = -1502.1428
$ws.Range("N82").Value = -3540

# LTW row 85
$ws.Range("H85").Value = 2424.8235
$ws.Range("I85").Value = 1863.1428
$ws.Range("J85").Value = 2818
$ws.Range("K85").Value = 1863.1428
$ws.Range("L85").Value = 2818
$ws.Range("M85").Value = -615.1428000000001
$ws.Range("N85").Value = -5314

# LTW row 132
$ws.Range("H132").Value = 3651.56
$ws.Range("I132").Value = 2707.4167
$ws.Range("J132").Value = 4523.077
$ws.Range("K132").Value = 8122.250100000001
$ws.Range("L132").Value = 13569.231
$ws.Range("M132").Value = -5592.250100000001
$ws.Range("N132").Value = -18629.231

# LTW row 140
$ws.Range("H140").Value = 36934.668
$ws.Range("J140").Value = 36934.668
$ws.Range("L140").Value = 36934.668
$ws.Range("N140").Value = -47294.668

$ws = $wb.Worksheets.Item("WVR")
# WVR row 15
$ws.Range("H15").Value = 29334.5
$ws.Range("J15").Value = 29334.5
$ws.Range("L15").Value = 29334.5
$ws.Range("N15").Value = -29910.5

# WVR row 31
$ws.Range("H31").Value = 43346
$ws.Range("J31").Value = 43346
$ws.Range("L31").Value = 43346
$ws.Range("N31").Value = -44042
